# Interdiff between v1 and v2 of UndoRedoActivityDiagram.pptx
# Re-lays-out the "undo stack" branch of the activity diagram:
#   - shifts the left-hand cluster (Oval/Arrow/RoundRect/Arrow) to the right
#   - removes the old "[undo or redo]" / "[command is undoable]" branch
#     (Diamond 11, TextBox 18, TextBox 19, Rectangle 21, Diamond 25 and
#     their 4 connecting elbow connectors)
#   - re-positions/resizes the remaining shapes of the right-hand branch
#     and updates their wording to talk about the "address book state"
#     instead of "command"

$EMU_PER_POINT = 12700.0

function Get-ShapeById {
    param($Shapes, [int]$Id)
    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $candidate = $Shapes.Item($i)
        if ($candidate.Id -eq $Id) {
            return $candidate
        }
    }
    throw "Shape with id $Id not found"
}

function Set-ShapeOffset {
    param($Shape, [double]$XEmu, [double]$YEmu)
    $Shape.Left = $XEmu / $EMU_PER_POINT
    $Shape.Top = $YEmu / $EMU_PER_POINT
}

function Set-ShapeExtent {
    param($Shape, [double]$CxEmu, [double]$CyEmu)
    $Shape.Width = $CxEmu / $EMU_PER_POINT
    $Shape.Height = $CyEmu / $EMU_PER_POINT
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shapes = $s.Shapes

# ---------------------------------------------------------------------
# 1) Move the left-hand cluster to the right
# ---------------------------------------------------------------------
Set-ShapeOffset (Get-ShapeById $shapes 4)  1929588 3261938   # Oval 3
Set-ShapeOffset (Get-ShapeById $shapes 6)  2165257 3379773   # Straight Arrow Connector 5
Set-ShapeOffset (Get-ShapeById $shapes 8)  2392862 3022393   # Rectangle: Rounded Corners 7
Set-ShapeOffset (Get-ShapeById $shapes 9)  3963217 3379774   # Straight Arrow Connector 8

# ---------------------------------------------------------------------
# 2) Remove the old "[undo or redo]" branch entirely: the diamond, the
#    two guard-condition textboxes, the "Clear redo stack" rounded
#    rectangle, the second diamond, and the four elbow connectors that
#    wired them together.
# ---------------------------------------------------------------------
(Get-ShapeById $shapes 12).Delete()   # Diamond 11
(Get-ShapeById $shapes 19).Delete()   # TextBox 18  "[undo or redo]"
(Get-ShapeById $shapes 20).Delete()   # TextBox 19  "[else]"
(Get-ShapeById $shapes 22).Delete()   # Rectangle: Rounded Corners 21  "Clear redo stack"
(Get-ShapeById $shapes 26).Delete()   # Diamond 25
(Get-ShapeById $shapes 24).Delete()   # Elbow Connector 23
(Get-ShapeById $shapes 55).Delete()   # Elbow Connector 54
(Get-ShapeById $shapes 57).Delete()   # Elbow Connector 56
(Get-ShapeById $shapes 61).Delete()   # Elbow Connector 60

# ---------------------------------------------------------------------
# 3) Re-position (and in some cases resize) the shapes of the right-hand
#    branch so it now starts where the removed branch used to be.
# ---------------------------------------------------------------------
Set-ShapeOffset (Get-ShapeById $shapes 46) 8215441 3248329   # Diamond 45

Set-ShapeOffset (Get-ShapeById $shapes 47) 4801950 3488712   # TextBox 46  "[else]"

$textBox47 = Get-ShapeById $shapes 48        # TextBox 47 "[command is undoable]"
Set-ShapeOffset $textBox47 2982124 2108748
Set-ShapeExtent $textBox47 2406969 646587
$textBox47.TextFrame.TextRange.Text = "[address book different from top of undo stack]"
Set-ShapeExtent $textBox47 2406969 646587

$roundRect50 = Get-ShapeById $shapes 51      # Rectangle: Rounded Corners 50 "Add command to undo stack"
Set-ShapeOffset $roundRect50 5389094 2377167
Set-ShapeExtent $roundRect50 2406970 888617
$roundRect50.TextFrame.TextRange.Text = "Add address book state to undo stack, clear redo stack"
Set-ShapeExtent $roundRect50 2406970 888617

Set-ShapeOffset (Get-ShapeById $shapes 56) 4356331 3140229   # Diamond 55

Set-ShapeOffset (Get-ShapeById $shapes 69) 8696207 3488712   # Straight Arrow Connector 68

Set-ShapeOffset (Get-ShapeById $shapes 75) 9115584 3378406   # Group 74

$elbow65 = Get-ShapeById $shapes 66          # Elbow Connector 65
Set-ShapeOffset $elbow65 4833528 2584663
Set-ShapeExtent $elbow65 318753 792380

$elbow71 = Get-ShapeById $shapes 72          # Elbow Connector 71
$elbow71.Rotation = 270
$elbow71.VerticalFlip = $false
Set-ShapeOffset $elbow71 6472219 1745490
Set-ShapeExtent $elbow71 108100 3859110
$elbow71.Adjustments.Item(1) = 3.11471

$elbow73 = Get-ShapeById $shapes 74          # Elbow Connector 73
Set-ShapeOffset $elbow73 7796064 2821476
Set-ShapeExtent $elbow73 659760 426853
